$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    -0.108838172,
    0.102213206,
    0.035056553,
    -0.07317942199999999,
    -0.075632767,
    -0.048480942,
    -0.026082961,
    -0.104705686,
    0.133157355,
    -0.082395143,
    0.268702188,
    -0.048298507,
    -0.192681301,
    -0.152074051,
    -0.052721073,
    0.183200194,
    -0.18698998,
    -0.138461451,
    -0.00775468,
    -0.010484921,
    0.124480542,
    -0.020797776,
    -0.009224501,
    0.043361183,
    -0.076622804,
    -0.32262107,
    -0.09768600199999999,
    -0.079447407,
    0.027373529,
    -0.079167162,
    -0.023174303,
    0.037933061,
    -0.183819699,
    -0.07962683500000001,
    -0.006608619,
    0.076644039,
    -0.002825303,
    -0.048441424,
    0.126694449,
    -0.026727514,
    -0.20244848,
    0.038890875,
    0.038352084,
    0.239657096,
    0.190517519,
    0.041690889,
    0.033719902,
    -0.188193707,
    0.096049495,
    -0.141508164,
    0.032553968,
    0.130092783,
    0.061158116,
    0.056411231,
    -0.02537516,
    -0.119164277,
    0.02497834,
    0.08360324199999999,
    -0.147494283,
    0.007501006,
    0.079947359,
    -0.07790449000000001,
    -0.050433668,
    -0.115189413,
    0.222678571,
    0.054921989,
    -0.119217342,
    -0.202754533,
    0.142159313,
    -0.135585723,
    -0.114476329,
    0.067702974,
    -0.154675907,
    -0.140164361,
    -0.323887601,
    0.04418762,
    0.400949576,
    0.071259618,
    -0.189579175,
    0.06696640299999999,
    -0.082723792,
    0.011222476,
    0.112965522,
    0.159667841,
    -0.009340342,
    0.010392343,
    -0.106651159,
    -0.011075524,
    0.194999854,
    -0.06378718799999999,
    -0.045034885,
    0.173100504,
    -0.010282529,
    0.094190442,
    0.044694176,
    0.052367868,
    -0.058467018,
    0.067675821,
    -0.113830893,
    -0.030414354,
    0.042017394,
    -0.06637910299999999,
    -0.026684042,
    0.119032121,
    -0.117996948,
    0.138562426,
    0.006768323,
    0.062455424,
    0.030197182,
    -0.010535627,
    -0.067032019,
    -0.06631896900000001,
    0.120439429,
    -0.206220964,
    0.226510605,
    0.178263064,
    0.087851041,
    0.109401843,
    0.12891979,
    0.136094716,
    -0.044693919,
    -0.021595968,
    -0.21338353,
    0.015409064,
    0.08630738,
    -0.031983804,
    0.08766318400000001,
    0.003020476
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 2, 2).Value = $values[$i]
}
